# Weekly update: insert a new price record at row 22 (Jengibre, Mercado
# Mayorista Lo Valledor de Santiago) and shift the existing history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22 - everything below (old rows
# 22-43) shifts down to 23-44, extending the used range to A1:R44.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with this week's record.
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C22").Value = "Metropolitana"
$ws.Range("D22").Value = 44438
$ws.Range("E22").Value = 13
$ws.Range("F22").Value = 100114007
$ws.Range("G22").Value = "Jengibre"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 13000
$ws.Range("L22").Value = 14000
$ws.Range("M22").Value = 13333
$ws.Range("N22").Value = "$/caja 13 kilos"
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 1026
$ws.Range("Q22").Value = 13
$ws.Range("R22").Value = "Hortaliza"
